# CachingApplication.xlsx - "Updated Caching Experiment"
#
# Fills in the second set of raw measurements (rows 33-42) on the
# WebViewMeasurement sheet, which lets the pre-existing AVERAGE() formulas
# in row 43 resolve instead of showing #DIV/0!, adds the corresponding
# series to the second ("Caching") bar chart, and updates a couple of
# view-state bits (zoom / selection) to match the author's session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WebViewMeasurement")

# ---------------------------------------------------------------------
# 1. Raw data for the second experiment table (rows 33-42)
# ---------------------------------------------------------------------
# Electron Proxy (columns B/C) - all 10 iterations
$ws.Range("B33").Value = 16.59
$ws.Range("C33").Value = 4.25
$ws.Range("B34").Value = 12.44
$ws.Range("C34").Value = 4.13
$ws.Range("B35").Value = 23.88
$ws.Range("C35").Value = 3.98
$ws.Range("B36").Value = 23.64
$ws.Range("C36").Value = 4.17
$ws.Range("B37").Value = 19.12
$ws.Range("C37").Value = 4.18
$ws.Range("B38").Value = 20.45
$ws.Range("C38").Value = 6.94
$ws.Range("B39").Value = 16.21
$ws.Range("C39").Value = 4.1900000000000004
$ws.Range("B40").Value = 20.46
$ws.Range("C40").Value = 4.08
$ws.Range("B41").Value = 16.41
$ws.Range("C41").Value = 7.1
$ws.Range("B42").Value = 15.99
$ws.Range("C42").Value = 7.06

# No Proxy (columns E/F) - only the first 7 iterations were recorded
$ws.Range("E33").Value = 11.56
$ws.Range("F33").Value = 11.35
$ws.Range("E34").Value = 9.7899999999999991
$ws.Range("F34").Value = 8.9600000000000009
$ws.Range("E35").Value = 9.6199999999999992
$ws.Range("F35").Value = 13.27
$ws.Range("E36").Value = 9.43
$ws.Range("F36").Value = 9.3000000000000007
$ws.Range("E37").Value = 10.02
$ws.Range("F37").Value = 11.38
$ws.Range("E38").Value = 9.5
$ws.Range("F38").Value = 9.64
$ws.Range("E39").Value = 9.0299999999999994
$ws.Range("F39").Value = 8.6999999999999993

# ---------------------------------------------------------------------
# 2. View state - zoom + selection on the WebViewMeasurement sheet
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("F40").Select()

# ---------------------------------------------------------------------
# 3. Add the 4 new series to the second bar chart (rows 31-43 table)
#    Colours mirror the theme accent1-4 colours used by the sibling
#    chart (Chart 3 / chart1.xml) for the same four measurements.
# ---------------------------------------------------------------------
$chart = $ws.ChartObjects("Chart 2").Chart

$ser1 = $chart.SeriesCollection().NewSeries()
$ser1.Name = "Electron Proxy First Request"
$ser1.Values = "=WebViewMeasurement!`$B`$43"
$ser1.Format.Fill.ForeColor.RGB = 13998939  # accent1 5B9BD5
$ser1.InvertIfNegative = $false

$ser2 = $chart.SeriesCollection().NewSeries()
$ser2.Name = "Electron Proxy Second Request"
$ser2.Values = "=WebViewMeasurement!`$C`$43"
$ser2.Format.Fill.ForeColor.RGB = 3243501   # accent2 ED7D31
$ser2.InvertIfNegative = $false

$ser3 = $chart.SeriesCollection().NewSeries()
$ser3.Name = "No Proxy First Request"
$ser3.Values = "=WebViewMeasurement!`$E`$43"
$ser3.Format.Fill.ForeColor.RGB = 10855845  # accent3 A5A5A5
$ser3.InvertIfNegative = $false

$ser4 = $chart.SeriesCollection().NewSeries()
$ser4.Name = "No Proxy Second Request"
$ser4.Values = "=WebViewMeasurement!`$F`$43"
$ser4.Format.Fill.ForeColor.RGB = 49407     # accent4 FFC000
$ser4.InvertIfNegative = $false

# Match the other chart's line/series styling: no outline on the bars.
$ser1.Format.Line.Visible = $false
$ser2.Format.Line.Visible = $false
$ser3.Format.Line.Visible = $false
$ser4.Format.Line.Visible = $false
